$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H26").Value = 0.82794
$ws.Range("I26").Value = 0.00716
$ws.Range("H27").Value = 0.36818
$ws.Range("I27").Value = 0.02333
$ws.Range("H28").Value = 0.83023
$ws.Range("I28").Value = 0.00777
$ws.Range("H29").Value = 0.4177
$ws.Range("I29").Value = 0.0243
$ws.Range("H30").Value = 0.82802
$ws.Range("I30").Value = 0.00697
$ws.Range("H31").Value = 0.3676
$ws.Range("I31").Value = 0.02275
$ws.Range("H32").Value = 0.83079
$ws.Range("I32").Value = 0.00778
$ws.Range("H33").Value = 0.42019
$ws.Range("I33").Value = 0.02377
$ws.Range("H34").Value = 0.78583
$ws.Range("I34").Value = 0.00611
$ws.Range("H35").Value = 0.19685
$ws.Range("I35").Value = 0.01948
$ws.Range("H36").Value = 0.78638
$ws.Range("I36").Value = 0.00607
$ws.Range("H37").Value = 0.20366
$ws.Range("I37").Value = 0.01931
$ws.Range("H38").Value = 0.78582
$ws.Range("I38").Value = 0.00612
$ws.Range("H39").Value = 0.19685
$ws.Range("I39").Value = 0.01948
$ws.Range("H40").Value = 0.78633
$ws.Range("I40").Value = 0.00604
$ws.Range("H41").Value = 0.20358
$ws.Range("I41").Value = 0.0192
$ws.Range("H66").Value = 0.79574
$ws.Range("I66").Value = 0.009719999999999999
$ws.Range("H67").Value = 0.23332
$ws.Range("I67").Value = 0.01805
$ws.Range("H68").Value = 0.80536
$ws.Range("I68").Value = 0.00978
$ws.Range("H69").Value = 0.34075
$ws.Range("I69").Value = 0.02601
$ws.Range("H70").Value = 0.79611
$ws.Range("I70").Value = 0.00963
$ws.Range("H71").Value = 0.23332
$ws.Range("I71").Value = 0.01808
$ws.Range("H72").Value = 0.8089
$ws.Range("I72").Value = 0.01002
$ws.Range("H73").Value = 0.35438
$ws.Range("I73").Value = 0.02522
$ws.Range("H74").Value = 0.77574
$ws.Range("I74").Value = 0.00787
$ws.Range("H75").Value = 0.14342
$ws.Range("I75").Value = 0.02251
$ws.Range("H76").Value = 0.77847
$ws.Range("I76").Value = 0.00763
$ws.Range("H77").Value = 0.17774
$ws.Range("I77").Value = 0.0219
$ws.Range("H78").Value = 0.77578
$ws.Range("I78").Value = 0.00788
$ws.Range("H79").Value = 0.14342
$ws.Range("I79").Value = 0.02251
$ws.Range("H80").Value = 0.77872
$ws.Range("I80").Value = 0.00762
$ws.Range("H81").Value = 0.17832
$ws.Range("I81").Value = 0.022
$ws.Range("H106").Value = 0.82436
$ws.Range("I106").Value = 0.0077
$ws.Range("H107").Value = 0.35763
$ws.Range("I107").Value = 0.0235
$ws.Range("H108").Value = 0.82658
$ws.Range("I108").Value = 0.0077
$ws.Range("H109").Value = 0.39842
$ws.Range("I109").Value = 0.02239
$ws.Range("H110").Value = 0.82461
$ws.Range("I110").Value = 0.00779
$ws.Range("H111").Value = 0.35879
$ws.Range("I111").Value = 0.02272
$ws.Range("H112").Value = 0.82698
$ws.Range("I112").Value = 0.00797
$ws.Range("H113").Value = 0.39992
$ws.Range("I113").Value = 0.02328
$ws.Range("H114").Value = 0.77855
$ws.Range("I114").Value = 0.00627
$ws.Range("H115").Value = 0.16909
$ws.Range("I115").Value = 0.02006
$ws.Range("H116").Value = 0.7786999999999999
$ws.Range("I116").Value = 0.00631
$ws.Range("H117").Value = 0.1715
$ws.Range("I117").Value = 0.02019
$ws.Range("H118").Value = 0.77857
$ws.Range("I118").Value = 0.00625
$ws.Range("H119").Value = 0.16909
$ws.Range("I119").Value = 0.02006
$ws.Range("H120").Value = 0.77873
$ws.Range("I120").Value = 0.00629
$ws.Range("H121").Value = 0.17167
$ws.Range("I121").Value = 0.02012
$ws.Range("H146").Value = 0.80875
$ws.Range("I146").Value = 0.00725
$ws.Range("H147").Value = 0.26997
$ws.Range("I147").Value = 0.01943
$ws.Range("H148").Value = 0.81526
$ws.Range("I148").Value = 0.00886
$ws.Range("H149").Value = 0.36951
$ws.Range("I149").Value = 0.02875
$ws.Range("H150").Value = 0.8089
$ws.Range("I150").Value = 0.00701
$ws.Range("H151").Value = 0.26997
$ws.Range("I151").Value = 0.01961
$ws.Range("H152").Value = 0.81912
$ws.Range("I152").Value = 0.009209999999999999
$ws.Range("H153").Value = 0.38538
$ws.Range("I153").Value = 0.03063
$ws.Range("H154").Value = 0.78178
$ws.Range("I154").Value = 0.008109999999999999
$ws.Range("H155").Value = 0.15073
$ws.Range("I155").Value = 0.02086
$ws.Range("H156").Value = 0.78347
$ws.Range("I156").Value = 0.008370000000000001
$ws.Range("H157").Value = 0.17425
$ws.Range("I157").Value = 0.02438
$ws.Range("H158").Value = 0.7819700000000001
$ws.Range("I158").Value = 0.008030000000000001
$ws.Range("H159").Value = 0.15073
$ws.Range("I159").Value = 0.02086
$ws.Range("H160").Value = 0.78377
$ws.Range("I160").Value = 0.008460000000000001
$ws.Range("H161").Value = 0.17441
$ws.Range("I161").Value = 0.02438
